$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels:
#   A1 was "KODE"   -> becomes "KODE ( Unik => Maks 4 karakter )"
#   B1 was "NEGARA" -> stays "NEGARA" (shared-string slot reordered in the source diff,
#                       but the visible cell value is unchanged)
$ws.Range("A1").Value = "KODE ( Unik => Maks 4 karakter )"
$ws.Range("B1").Value = "NEGARA"

# Widen the columns to fit the new header text
$ws.Columns.Item(1).ColumnWidth = 31.6667
$ws.Columns.Item(2).ColumnWidth = 30.5

# Move the active selection to D9, matching the saved view state
$ws.Range("D9").Select()
